# Update the "dSF" column (F) values for several rows, per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value  = -3
$ws.Range("F6").Value  = -8
$ws.Range("F7").Value  = -5
$ws.Range("F9").Value  = 3
$ws.Range("F10").Value = -1
$ws.Range("F12").Value = 6
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = -2
$ws.Range("F17").Value = -5
